$wb = $excel.ActiveWorkbook

# --- RealData sheet: append two new daily data points ---
$real = $wb.Worksheets.Item("RealData")
$real.Range("A112").Value = 44128
$real.Range("B112").Value = 1128
$real.Range("A113").Value = 44129
$real.Range("B113").Value = 1208

# --- Model sheet: lower the growth-rate input (column C) for rows 105-120 ---
$model = $wb.Worksheets.Item("Model")

# Rows 114-120 already carry the "projected" highlight style (s="4"); the
# edit flattens that formatting on column C (same effect as typing the new
# value over the existing highlighted cells / copying the unstyled C105
# value down), while rows 105-113 (already unstyled) are unaffected.
$model.Range("C114:C120").ClearFormats()
$model.Range("C105:C120").Value = 1.95

# --- restore on-screen selections to match the edited areas ---
$real.Activate()
$real.Range("F107").Select()

$model.Activate()
$model.Range("C105:C120").Select()
